$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Solothurn (row 12): delivery interval text gets the "sofern neue
# Gemeindeergebnisse vorhanden sind" / "si de nouveaux resultats communales
# sont disponibles" qualifier appended, same as already used for Vaud
# (row 23). Copy that row's cell formatting first (longer text needs the
# wrapping / bottom-border style), then update the text itself.
$ws.Range("C23").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("D23").Copy()
$ws.Range("D12").PasteSpecial(-4122)

$ws.Range("C12").Value2 = "Lieferintervall: alle 5 Minuten, sofern neue Gemeindeergebnisse vorhanden sind"
$ws.Range("D12").Value2 = "Intervalle de transfert: toutes les 5 minutes, si de nouveaux résultats communales sont disponibles"

# Row grew to two lines of wrapped text.
$ws.Rows(12).RowHeight = 29.25

# Selection moved to D12.
$ws.Range("D12").Select()
